# Slutprojekt.docx edit script
# Applies:
#  1) Merge the two runs split by the _GoBack bookmark around
#     ", vilket skulle kunna orsaka problem om någon..." (text unchanged,
#     bookmark removed as a side effect of the text replace).
#  2) Merge "sig som kund" + "er" into "sig som kunder" (text unchanged,
#     bookmark removed as a side effect of the text replace).
#  3) Add a new paragraph about missing payment-structure design after the
#     "...dessa två roller. " paragraph.
#  4) Fill in the previously empty paragraph under the "Storage Engine"
#     heading with the MyISAM/InnoDB discussion.
#  5) Bump the cached PAGE field result in the footer from 3 to 4.

$d = $word.ActiveDocument

# --- 1) & 2): collapse the bookmark-split runs back into single runs ---
$d.Content.Find.Execute(
    ", vilket skulle kunna orsaka problem om någon sparad data skulle ändras och behöva uppdateras",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ", vilket skulle kunna orsaka problem om någon sparad data skulle ändras och behöva uppdateras",
    2) | Out-Null

$d.Content.Find.Execute(
    "sig som kunder",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "sig som kunder",
    2) | Out-Null

# --- helper: apply the document's body font to a range --------------------
function Set-SegoeFont($rng) {
    $rng.Font.NameAscii = "Segoe UI"
    $rng.Font.Name = "Segoe UI"
    $rng.Font.NameBi = "Segoe UI"
}

# --- 3) new paragraph about payments, inserted after the paragraph that
#        ends in "...dessa två roller. " -----------------------------------

# Locate the "Storage Engine" heading paragraph, the paragraph right before
# it is the one ending in "...dessa två roller. ".
$storageHeadingIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Storage Engine")) {
        $storageHeadingIndex = $i
        break
    }
}

$prevPara = $d.Paragraphs($storageHeadingIndex - 1)
$prevRange = $prevPara.Range
$prevRange.Collapse(0)
$prevRange.InsertParagraphAfter()

$newParaIndex = $storageHeadingIndex
$newPara = $d.Paragraphs($newParaIndex)
$insertPos = $newPara.Range.Start
$insPoint = $d.Range($insertPos, $insertPos)

$paymentText = "Ännu en nackdel med min design och struktur är att jag inte implementerat någon struktur för betalning av uthyrningarna. Jag valde att inte ta med det i min databas, då det finns en så stor mängd sätt att genomföra betalningar, med kreditkort, fakturor, Swish, osv. och det kändes som ett stort tillägg och till viss del utanför uppgiftens scope att skapa en filmuthyrningsdatabas. Dock skulle det troligtvis vara relativt enkelt att implementera detta genom att skapa relevanta betalningstabeller, och sedan koppla in dem i de existerande tabellerna, och möjligtvis till varje uthyrning koppla en instans av en betalning. "

$insPoint.InsertAfter($paymentText)
Set-SegoeFont $insPoint

# --- 4) fill the empty paragraph under "Storage Engine" --------------------

# Re-find the heading (index may have shifted because of the insert above).
$storageHeadingIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Storage Engine")) {
        $storageHeadingIndex = $i
        break
    }
}

$emptyPara = $d.Paragraphs($storageHeadingIndex + 1)
$emptyPos = $emptyPara.Range.Start
$emptyPoint = $d.Range($emptyPos, $emptyPos)

$engineText = "Den storage engine jag har valt för alla mina tabeller är den som är standard i MySQL just nu, InnoDB. Från det jag har läst om MyISAM verkar det som att det är en snabbare storage engine än InnoDB när det kommer till läs- och skriv-operationer. Utöver det kan jag inte se några fördelar med MyISAM över InnoDB när det komemr till den databas jag har skapat. Om du vill göra en uppdatering i någon tabell i MyISAM låser den hela tabellen tills du är klar med uppdateringen, medan InnoDB endast låser de rader du kommer att påverka, vilket är fördelaktigt om t.ex. flera kunder försöker ändra på sina uppgifter samtidigt, då väntetiderna blir kortare. InnoDBs stöd för transaktioner som MyISAM inte har ser jag ingen större poäng med att använda så som databasen ser ut idag, då inga komplicerade transaktioner sker. Skulle den däremot utökas och implementera någon slags credit-system där man kan ladda på sitt konto med pengar och sedan använda dem för att hyra filmer med skulle det vara en bra idé att göra det med transaktioner. Den största anledningen varför InnoDB är det bästa valet för min databas är för att den stödjer integritet och restriktioner för främmande nycklar. Då min databas är uppdelad i många olika tabeller med mycket relationer mellan dem är det viktigt att de olika främmande nycklarna behåller sina referenser och att man kan välja vad som ska hända när en viss post ska tas bort från en tabell som refereras till i en annan tabell."

$emptyPoint.InsertAfter($engineText)
Set-SegoeFont $emptyPoint

# --- 5) bump the cached PAGE field number in the footer --------------------
foreach ($story in $d.StoryRanges) {
    if ($story.StoryType -eq 9) {
        foreach ($fld in $story.Fields) {
            if ($fld.Type -eq 33) {
                $resultRange = $fld.Result
                if ($resultRange.Text -eq "3") {
                    $resultRange.Characters(1).Text = "4"
                }
            }
        }
    }
}
